$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TURF_Summary_Export")

# Update the subgroup label text (shared string "Caregivers" -> "Current WWP Members")
$ws.Range("E18").Value = "Current WWP Members"

# Row 3: C3 was a string ("Offered") -> becomes the number 66
$ws.Range("C3").Value = 66

# Row 4: C4 was a number (79.17) -> becomes the string "Considered"
$ws.Range("C4").Value = "Considered"

# Row 8: C8 was a string ("Offered") -> becomes the number 71.59
$ws.Range("C8").Value = 71.59

# Row 11: C11 was a string ("Offered") -> becomes the number 46.69
$ws.Range("C11").Value = 46.69

# Row 14: C14 was a string ("Offered") -> becomes the number 58.71
$ws.Range("C14").Value = 58.71

# Row 18 updates
$ws.Range("B18").Value = 1.339
$ws.Range("C18").Value = 76.33
$ws.Range("D18").Value = 72.06
$ws.Range("F18").Value = 177
